# Updates fitting of MEG, IMP in chapter5.tex
# - Adds "Blazed:" / "Blaze angles:" header labels in L1:M1
# - Fills M3:M51 with blaze-angle values 0.2 .. 5.0 (step 0.1)
# - Updates the active selection / window view to reflect the new data range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels
$ws.Range("L1").Value = "Blazed:"
$ws.Range("M1").Value = "Blaze angles:"

# New "Blaze angles" column of values (0.2 through 5.0 in steps of 0.1)
$blazeAngles = @(0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1, 1.1, 1.2, 1.3, 1.4, 1.5, 1.6, 1.7, 1.8, 1.9, 2, 2.1, 2.2, 2.3, 2.4, 2.5, 2.6, 2.7, 2.8, 2.9, 3, 3.1, 3.2, 3.3, 3.4, 3.5, 3.6, 3.7, 3.8, 3.9, 4, 4.1, 4.2, 4.3, 4.4, 4.5, 4.6, 4.7, 4.8, 4.9, 5)

for ($i = 0; $i -lt $blazeAngles.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 13).Value = $blazeAngles[$i]
}

# Update selection to match the reviewed range of the new column
$ws.Range("D43:D141").Select()
